$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 2900.6667
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2900.6667
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2900.6667
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3552.6667

# Row 62
$ws.Range("H62").Value = 866.25
$ws.Range("J62").Value = 918.5714
$ws.Range("L62").Value = 918.5714
$ws.Range("N62").Value = -2166.5714

# Row 65
$ws.Range("H65").Value = 866.25
$ws.Range("J65").Value = 918.5714
$ws.Range("L65").Value = 4592.857
$ws.Range("N65").Value = -10832.857

# Row 112
$ws.Range("H112").Value = 9434849
$ws.Range("I112").Value = 643.3333
$ws.Range("K112").Value = 1929.9999
$ws.Range("M112").Value = -821.9999

# Row 129
$ws.Range("H129").Value = 1045.9578
$ws.Range("J129").Value = 1061.9851
$ws.Range("L129").Value = 3185.955300000001
$ws.Range("N129").Value = -13185.9553

# Row 132
$ws.Range("H132").Value = 2599.9688
$ws.Range("I132").Value = 1823.88
$ws.Range("J132").Value = 5371.7144
$ws.Range("K132").Value = 5471.64
$ws.Range("L132").Value = 16115.1432
$ws.Range("M132").Value = -2941.64
$ws.Range("N132").Value = -21175.1432

$ws = $wb.Worksheets.Item("ARM")
# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 10884.074
$ws.Range("I20").Value = 1443.25
$ws.Range("J20").Value = 24616.182
$ws.Range("K20").Value = 1443.25
$ws.Range("L20").Value = 24616.182
$ws.Range("M20").Value = -1196.25
$ws.Range("N20").Value = -25110.182

# Row 105
$ws.Range("H105").Value = 24708.777
$ws.Range("I105").Value = 68040.336
$ws.Range("J105").Value = 3043
$ws.Range("K105").Value = 68040.336
$ws.Range("L105").Value = 3043
$ws.Range("M105").Value = -66293.336
$ws.Range("N105").Value = -6537

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 190
$ws.Range("I7").Value = 165.55556
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 165.55556
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -52.55556000000001
$ws.Range("N7").Value = -526

# Row 58
$ws.Range("H58").Value = 2500
$ws.Range("I58").Value = 2500
$ws.Range("K58").Value = 2500
$ws.Range("M58").Value = -2297

# Row 132
$ws.Range("H132").Value = 4496
$ws.Range("I132").Value = 4992
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 14976
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -12446
$ws.Range("N132").Value = -17060

# Row 136
$ws.Range("H136").Value = 2500
$ws.Range("I136").Value = 2500
$ws.Range("K136").Value = 7500
$ws.Range("M136").Value = -4950

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 4123.3335
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 4348
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 13044
$ws.Range("M80").Value = -8064
$ws.Range("N80").Value = -14916

# Row 83
$ws.Range("H83").Value = 4123.3335
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 4348
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 39132
$ws.Range("M83").Value = -22320
$ws.Range("N83").Value = -48492

# Row 108
$ws.Range("H108").Value = 50362.5
$ws.Range("I108").Value = 50362.5
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 151087.5
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -148207.5
$ws.Range("N108").ClearContents()

# Row 111
$ws.Range("H111").Value = 700
$ws.Range("I111").Value = 700
$ws.Range("K111").Value = 2100
$ws.Range("M111").Value = 967

# Row 131
$ws.Range("H131").Value = 16924172
$ws.Range("I131").Value = 6667277
$ws.Range("J131").Value = 20001240
$ws.Range("K131").Value = 20001831
$ws.Range("L131").Value = 60003720
$ws.Range("M131").Value = -19996791
$ws.Range("N131").Value = -60013800

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5480.5
$ws.Range("I70").Value = 5565.36
$ws.Range("J70").Value = 5215.3125
$ws.Range("K70").Value = 5565.36
$ws.Range("L70").Value = 5215.3125
$ws.Range("M70").Value = -5295.36
$ws.Range("N70").Value = -5755.3125

# Row 73
$ws.Range("H73").Value = 5480.5
$ws.Range("I73").Value = 5565.36
$ws.Range("J73").Value = 5215.3125
$ws.Range("K73").Value = 5565.36
$ws.Range("L73").Value = 5215.3125
$ws.Range("M73").Value = -4629.36
$ws.Range("N73").Value = -7087.3125

# Row 97
$ws.Range("H97").Value = 640
$ws.Range("I97").Value = 640
$ws.Range("K97").Value = 640
$ws.Range("M97").Value = -144

# Row 102
$ws.Range("H102").Value = 1833.6666
$ws.Range("I102").Value = 1350.5
$ws.Range("K102").Value = 1350.5
$ws.Range("M102").Value = 271.5

# Row 132
$ws.Range("H132").Value = 3678.2727
$ws.Range("I132").Value = 3659.238
$ws.Range("J132").Value = 3695.652
$ws.Range("K132").Value = 10977.714
$ws.Range("L132").Value = 11086.956
$ws.Range("M132").Value = -8447.714
$ws.Range("N132").Value = -16146.956

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2586548
$ws.Range("J22").Value = 2854.054
$ws.Range("L22").Value = 2854.054
$ws.Range("N22").Value = -3444.054

# Row 27
$ws.Range("H27").Value = 2586548
$ws.Range("J27").Value = 2854.054
$ws.Range("L27").Value = 2854.054
$ws.Range("N27").Value = -3068.054

# Row 46
$ws.Range("H46").Value = 15874461
$ws.Range("I46").Value = 27778816
$ws.Range("K46").Value = 27778816
$ws.Range("M46").Value = -27778628

# Row 55
$ws.Range("H55").Value = 26316182
$ws.Range("I55").Value = 265.9
$ws.Range("J55").Value = 55556090
$ws.Range("K55").Value = 265.9
$ws.Range("L55").Value = 55556090
$ws.Range("M55").Value = -92.89999999999998
$ws.Range("N55").Value = -55556436

# Row 68
$ws.Range("H68").Value = 24391822
$ws.Range("I68").Value = 1625.32
$ws.Range("J68").Value = 62501504
$ws.Range("K68").Value = 1625.32
$ws.Range("L68").Value = 62501504
$ws.Range("M68").Value = -876.3199999999999
$ws.Range("N68").Value = -62503002

# Row 71
$ws.Range("H71").Value = 24391822
$ws.Range("I71").Value = 1625.32
$ws.Range("J71").Value = 62501504
$ws.Range("K71").Value = 8126.599999999999
$ws.Range("L71").Value = 312507520
$ws.Range("M71").Value = -4382.599999999999
$ws.Range("N71").Value = -312515008

# Row 82
$ws.Range("H82").Value = 4943.1333
$ws.Range("I82").Value = 1152.3334
$ws.Range("J82").Value = 10629.333
$ws.Range("K82").Value = 1152.3334
$ws.Range("L82").Value = 10629.333
$ws.Range("M82").Value = -791.3334
$ws.Range("N82").Value = -11351.333

# Row 85
$ws.Range("H85").Value = 4943.1333
$ws.Range("I85").Value = 1152.3334
$ws.Range("J85").Value = 10629.333
$ws.Range("K85").Value = 1152.3334
$ws.Range("L85").Value = 10629.333
$ws.Range("M85").Value = 95.66660000000002
$ws.Range("N85").Value = -13125.333

# Row 93
$ws.Range("H93").Value = 55580104
$ws.Range("I93").Value = 50600
$ws.Range("J93").Value = 100003710
$ws.Range("K93").Value = 50600
$ws.Range("L93").Value = 100003710
$ws.Range("M93").Value = -49352
$ws.Range("N93").Value = -100006206

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()

# Row 136
$ws.Range("H136").Value = 965.6667
$ws.Range("I136").Value = 909.6842
$ws.Range("J136").Value = 1497.5
$ws.Range("K136").Value = 2729.0526
$ws.Range("L136").Value = 4492.5
$ws.Range("M136").Value = -179.0526
$ws.Range("N136").Value = -9592.5
